$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-30 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-31 Friday", 2) | Out-Null
$d.Content.Find.Execute("801×2=1602", $true, $false, $false, $false, $false, $true, 1, $false, "846×4=3384", 2) | Out-Null
$d.Content.Find.Execute("458×6=2748", $true, $false, $false, $false, $false, $true, 1, $false, "172×9=1548", 2) | Out-Null
$d.Content.Find.Execute("559×3=1677", $true, $false, $false, $false, $false, $true, 1, $false, "398×4=1592", 2) | Out-Null
$d.Content.Find.Execute("888×4=3552", $true, $false, $false, $false, $false, $true, 1, $false, "419×7=2933", 2) | Out-Null
$d.Content.Find.Execute("567×5=2835", $true, $false, $false, $false, $false, $true, 1, $false, "613×7=4291", 2) | Out-Null
$d.Content.Find.Execute("191×2=382", $true, $false, $false, $false, $false, $true, 1, $false, "965×4=3860", 2) | Out-Null
$d.Content.Find.Execute("112×3=336", $true, $false, $false, $false, $false, $true, 1, $false, "393×5=1965", 2) | Out-Null
$d.Content.Find.Execute("252×6=1512", $true, $false, $false, $false, $false, $true, 1, $false, "179×9=1611", 2) | Out-Null
$d.Content.Find.Execute("377×9=3393", $true, $false, $false, $false, $false, $true, 1, $false, "941×6=5646", 2) | Out-Null
$d.Content.Find.Execute("138×3=414", $true, $false, $false, $false, $false, $true, 1, $false, "243×2=486", 2) | Out-Null
$d.Content.Find.Execute("227×3=681", $true, $false, $false, $false, $false, $true, 1, $false, "494×4=1976", 2) | Out-Null
$d.Content.Find.Execute("520×5=2600", $true, $false, $false, $false, $false, $true, 1, $false, "750×8=6000", 2) | Out-Null
$d.Content.Find.Execute("416×3=1248", $true, $false, $false, $false, $false, $true, 1, $false, "749×4=2996", 2) | Out-Null
$d.Content.Find.Execute("279×3=837", $true, $false, $false, $false, $false, $true, 1, $false, "621×5=3105", 2) | Out-Null
$d.Content.Find.Execute("835×7=5845", $true, $false, $false, $false, $false, $true, 1, $false, "750×2=1500", 2) | Out-Null
$d.Content.Find.Execute("480×4=1920", $true, $false, $false, $false, $false, $true, 1, $false, "416×2=832", 2) | Out-Null
$d.Content.Find.Execute("233×8=1864", $true, $false, $false, $false, $false, $true, 1, $false, "242×4=968", 2) | Out-Null
$d.Content.Find.Execute("485×7=3395", $true, $false, $false, $false, $false, $true, 1, $false, "580×2=1160", 2) | Out-Null
$d.Content.Find.Execute("261×5=1305", $true, $false, $false, $false, $false, $true, 1, $false, "584×2=1168", 2) | Out-Null
$d.Content.Find.Execute("617×3=1851", $true, $false, $false, $false, $false, $true, 1, $false, "782×6=4692", 2) | Out-Null
$d.Content.Find.Execute("672×4=2688", $true, $false, $false, $false, $false, $true, 1, $false, "816×4=3264", 2) | Out-Null
$d.Content.Find.Execute("431×6=2586", $true, $false, $false, $false, $false, $true, 1, $false, "131×2=262", 2) | Out-Null
$d.Content.Find.Execute("226×6=1356", $true, $false, $false, $false, $false, $true, 1, $false, "748×6=4488", 2) | Out-Null
$d.Content.Find.Execute("410×5=2050", $true, $false, $false, $false, $false, $true, 1, $false, "709×4=2836", 2) | Out-Null
$d.Content.Find.Execute("277×7=1939", $true, $false, $false, $false, $false, $true, 1, $false, "787×7=5509", 2) | Out-Null
